# LF Energy High Level Overview Deck - April 2025 update (#56)
#
# 1) Re-point every table in the deck from the old table style GUID to the
#    new one.
# 2) Swap the "Default" / "Simple Light" colour schemes that live in the
#    deck's themes (the deck's active theme - theme3.xml, reached through
#    the slide masters / ThemeColorScheme - takes on the values that used
#    to belong to the "Default" scheme).

$p = $ppt.ActivePresentation

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# ---------------------------------------------------------------------
# 1. Table styles: {F126362D-DFA5-4DE9-A93B-187AFD16EC08} -> {1B8F7AFB-DB68-486D-90CB-5E2D6C7F809D}
# ---------------------------------------------------------------------
$oldStyleId = "{F126362D-DFA5-4DE9-A93B-187AFD16EC08}"
$newStyleId = "{1B8F7AFB-DB68-486D-90CB-5E2D6C7F809D}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Theme colours: the deck's main colour scheme ("Simple Light") takes
#    on the colour values that used to be the "Default" scheme.
#    (dk1/lt1 are identical in both schemes, so only dk2/lt2/accent1-6/
#    hlink/folHlink actually change.)
# ---------------------------------------------------------------------
$newColors = @{
    3  = "158158"  # dk2
    4  = "F3F3F3"  # lt2
    5  = "058DC7"  # accent1
    6  = "50B432"  # accent2
    7  = "ED561B"  # accent3
    8  = "EDEF00"  # accent4
    9  = "24CBE5"  # accent5
    10 = "64E572"  # accent6
    11 = "2200CC"  # hlink
    12 = "551A8B"  # folHlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $newColors.Keys) {
    $themeColors.Item($idx).RGB = HexToRgbInt($newColors[$idx])
}
